# Update submgmt test data: status code 400 -> 200

$wb = $excel.ActiveWorkbook

# --- Sheet: registerSubscriptions (F7:F10 -> 400 to 200) ---
$ws1 = $wb.Worksheets.Item("registerSubscriptions")
$ws1.Range("F7").Value = 200
$ws1.Range("F8").Value = 200
$ws1.Range("F9").Value = 200
$ws1.Range("F10").Value = 200

# These rows also drop the explicit fill-apply flag on F so the cell's
# style matches the neighboring (non-filled) cells (xf 5 -> xf 4).
$ws1.Range("F7:F10").Interior.Pattern = -4142

# --- Sheet: delSubscriptionById (G2, G3, G6 -> 400 to 200) ---
$ws2 = $wb.Worksheets.Item("delSubscriptionById")
$ws2.Range("G2").Value = 200
$ws2.Range("G3").Value = 200
$ws2.Range("G6").Value = 200

# --- Active sheet / selection housekeeping ---
# The workbook now opens on "delSubscriptionById" first with its own
# selection, then settles back on "registerSubscriptions" as the tab
# that is actually showing (tabSelected) with a new selected cell.
$ws2.Activate() | Out-Null
$ws2.Range("I6").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("F3").Select() | Out-Null
